$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-05 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-06 Saturday", 2) | Out-Null
$d.Content.Find.Execute("619÷7=88, 3", $true, $false, $false, $false, $false, $true, 1, $false, "633÷7=90, 3", 2) | Out-Null
$d.Content.Find.Execute("291÷6=48, 3", $true, $false, $false, $false, $false, $true, 1, $false, "250÷6=41, 4", 2) | Out-Null
$d.Content.Find.Execute("957÷6=159, 3", $true, $false, $false, $false, $false, $true, 1, $false, "407÷4=101, 3", 2) | Out-Null
$d.Content.Find.Execute("171÷7=24, 3", $true, $false, $false, $false, $false, $true, 1, $false, "523÷4=130, 3", 2) | Out-Null
$d.Content.Find.Execute("405÷9=45, 0", $true, $false, $false, $false, $false, $true, 1, $false, "885÷2=442, 1", 2) | Out-Null
$d.Content.Find.Execute("776÷8=97, 0", $true, $false, $false, $false, $false, $true, 1, $false, "762÷5=152, 2", 2) | Out-Null
$d.Content.Find.Execute("875÷4=218, 3", $true, $false, $false, $false, $false, $true, 1, $false, "972÷4=243, 0", 2) | Out-Null
$d.Content.Find.Execute("161÷6=26, 5", $true, $false, $false, $false, $false, $true, 1, $false, "528÷6=88, 0", 2) | Out-Null
$d.Content.Find.Execute("541÷3=180, 1", $true, $false, $false, $false, $false, $true, 1, $false, "493÷7=70, 3", 2) | Out-Null
$d.Content.Find.Execute("724÷2=362, 0", $true, $false, $false, $false, $false, $true, 1, $false, "134÷8=16, 6", 2) | Out-Null
$d.Content.Find.Execute("489÷2=244, 1", $true, $false, $false, $false, $false, $true, 1, $false, "475÷8=59, 3", 2) | Out-Null
$d.Content.Find.Execute("584÷5=116, 4", $true, $false, $false, $false, $false, $true, 1, $false, "625÷9=69, 4", 2) | Out-Null
$d.Content.Find.Execute("983÷6=163, 5", $true, $false, $false, $false, $false, $true, 1, $false, "458÷6=76, 2", 2) | Out-Null
$d.Content.Find.Execute("452÷2=226, 0", $true, $false, $false, $false, $false, $true, 1, $false, "984÷9=109, 3", 2) | Out-Null
$d.Content.Find.Execute("618÷9=68, 6", $true, $false, $false, $false, $false, $true, 1, $false, "568÷5=113, 3", 2) | Out-Null
$d.Content.Find.Execute("142÷6=23, 4", $true, $false, $false, $false, $false, $true, 1, $false, "345÷2=172, 1", 2) | Out-Null
$d.Content.Find.Execute("708÷5=141, 3", $true, $false, $false, $false, $false, $true, 1, $false, "387÷9=43, 0", 2) | Out-Null
$d.Content.Find.Execute("808÷2=404, 0", $true, $false, $false, $false, $false, $true, 1, $false, "900÷7=128, 4", 2) | Out-Null
$d.Content.Find.Execute("379÷4=94, 3", $true, $false, $false, $false, $false, $true, 1, $false, "858÷4=214, 2", 2) | Out-Null
$d.Content.Find.Execute("582÷5=116, 2", $true, $false, $false, $false, $false, $true, 1, $false, "512÷6=85, 2", 2) | Out-Null
$d.Content.Find.Execute("456÷4=114, 0", $true, $false, $false, $false, $false, $true, 1, $false, "268÷2=134, 0", 2) | Out-Null
$d.Content.Find.Execute("307÷7=43, 6", $true, $false, $false, $false, $false, $true, 1, $false, "296÷2=148, 0", 2) | Out-Null
$d.Content.Find.Execute("854÷5=170, 4", $true, $false, $false, $false, $false, $true, 1, $false, "623÷3=207, 2", 2) | Out-Null
$d.Content.Find.Execute("900÷4=225, 0", $true, $false, $false, $false, $false, $true, 1, $false, "389÷3=129, 2", 2) | Out-Null
$d.Content.Find.Execute("183÷5=36, 3", $true, $false, $false, $false, $false, $true, 1, $false, "713÷6=118, 5", 2) | Out-Null
